$d = $word.ActiveDocument

# 1. Add a first-line indent to the first paragraph ("Git commit --amend ...")
$p1 = $d.Paragraphs.Item(1)
$p1.FirstLineIndent = 21
$p1.CharacterUnitFirstLineIndent = 0

# 2. Move the "_GoBack" bookmark from the last paragraph that has it to the
#    very start of the document (start of the first paragraph), as a
#    zero-length bookmark.
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
}

# Work around zero-length-range-at-document-start quirk: insert a throwaway
# character at position 0, anchor the bookmark right after it (position 1),
# then remove the throwaway character so the bookmark collapses cleanly to
# the very beginning of the document while keeping Start/End adjacent.
$rStart = $d.Range(0, 0)
$rStart.InsertBefore("X")

$rBookmark = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $rBookmark)

$rTemp = $d.Range(0, 1)
$rTemp.Delete()
